# APACHE POI  SELENIUM IMPLEMENTATION
#
# Records three successive executions of the Gas-Mileage data-driven test:
#   - row 2 already held the 1st pass' input columns (A:D); this fills in
#     its computed Expected/Actual/Status/Time (E:H) -- only the Time
#     actually differs from what was there before.
#   - rows 3-5 already held the 1st pass' remaining input columns (A:D);
#     this completes their Expected/Actual/Status/Time (E:H).
#   - rows 6-13 are two more full passes (4 data rows each) appended below,
#     identical input data, each row stamped with its own execution time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 4-row input/result pattern produced by every pass of the test.
$pattern = @(
    @{ A = "Y"; B = 123000; C = 122000; D = 75;  EF = "13.33" },
    @{ A = "Y"; B = 15000;  C = 14000;  D = 60;  EF = "16.67" },
    @{ A = "Y"; B = 9877;   C = 7000;   D = 150; EF = "19.18" },
    @{ A = "Y"; B = 500;    C = 0;      D = 15;  EF = "33.33" }
)

# One execution timestamp per data row, rows 2 through 13 (3 full passes).
$timestamps = @(
    "03:26:08 PM", "03:26:09 PM", "03:26:11 PM", "03:26:12 PM",
    "03:26:13 PM", "03:26:14 PM", "03:26:17 PM", "03:26:18 PM",
    "03:26:20 PM", "03:26:21 PM", "03:26:22 PM", "03:26:23 PM"
)

# Columns E and F hold numeric-looking text ("13.33", ...) that must stay
# text (shared-string) rather than being auto-coerced to numbers -- format
# the whole block as Text up front, write the values, then restore the
# "Normal" style so no stray numeric formatting is left behind on the cells.
$efRange = $ws.Range("E2:F13")
$efRange.NumberFormat = "@"

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $p = $pattern[$i % 4]

    $ws.Cells.Item($row, 1).Value = $p.A
    $ws.Cells.Item($row, 2).Value = $p.B
    $ws.Cells.Item($row, 3).Value = $p.C
    $ws.Cells.Item($row, 4).Value = $p.D
    $ws.Cells.Item($row, 5).Value = $p.EF
    $ws.Cells.Item($row, 6).Value = $p.EF
    $ws.Cells.Item($row, 7).Value = "PASS!"
    $ws.Cells.Item($row, 8).Value = $timestamps[$i]
}

$efRange.Style = "Normal"

$ws.Range("F8").Select()
